# Applies the "apply a new Office-style theme to the deck" edit:
#   1) The three tables (on what are slides 14, 15 and 16) that referenced
#      the deck's custom default table style now use the built-in
#      "Medium Style 2 - Accent 1" table style instead.
#   2) The presentation's theme colour palette changes from the old
#      "Integral / Red Violet" scheme to the standard Office colour scheme.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Re-style the three tables that used the old default table style.
# ---------------------------------------------------------------------
$newTableStyleId = "{F5DFF213-0487-443E-A8AA-B33DC9FB411C}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------------
# 2) Swap the theme colour scheme over to the standard Office palette.
# ---------------------------------------------------------------------
$tcs = $p.Slides.Item(1).ThemeColorScheme

# Index : role      : new RGB (as a VBA-style BGR long = R + G*256 + B*65536)
$tcs.Item(1).RGB  = 0        # dk1      -> 000000
$tcs.Item(2).RGB  = 16777215 # lt1      -> FFFFFF
$tcs.Item(3).RGB  = 6968388  # dk2      -> 44546A
$tcs.Item(4).RGB  = 15132391 # lt2      -> E7E6E6
$tcs.Item(5).RGB  = 13998939 # accent1  -> 5B9BD5
$tcs.Item(6).RGB  = 3243501  # accent2  -> ED7D31
$tcs.Item(7).RGB  = 10855845 # accent3  -> A5A5A5
$tcs.Item(8).RGB  = 49407    # accent4  -> FFC000
$tcs.Item(9).RGB  = 12874308 # accent5  -> 4472C4
$tcs.Item(10).RGB = 4697456  # accent6  -> 70AD47
$tcs.Item(11).RGB = 12673797 # hlink    -> 0563C1
$tcs.Item(12).RGB = 7491477  # folHlink -> 954F72
